# Daily attendance processing - 2025-10-10 22:49:41
# Reorders the "Recorded By" (column G) audit list on the
# "Session Analysis Results" sheet so that the literal "System" marker is
# moved from the front of the comma-separated list to the back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "System, *") {
        $parts = $val -split ", "
        $reversed = $parts[($parts.Count - 1)..0]
        $cell.Value2 = $reversed -join ", "
    }
}
